# Add a new "2022-Q3" sheet (with its own fund-holding snapshot) right after
# the "总计" (totals) sheet, and record its summary row on the "总计" sheet.
#
# Existing "2022-Q2" / "2021-Q3" / "2021-Q2" / "2020-Q4" sheets are left
# completely untouched - inserting the new sheet simply pushes them one
# position to the right in the tab strip.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)          # "总计"
$q2202 = $wb.Worksheets.Item(2)          # "2022-Q2" (template for the new sheet's layout/styles)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by copying "2022-Q2" (so it inherits the
#    same column widths / header style / border style / page margins) and
#    re-position it immediately after "总计".
# ---------------------------------------------------------------------
$q2202.Copy($null, $total)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template sheet has two fund rows; 2022-Q3 only has one, so drop row 3.
$q3.Rows.Item(3).Delete()

# Force the numeric-looking text columns (B:G) to stay text, matching the
# rest of the workbook where these fund figures are stored as strings.
$q3.Range("B2:G2").NumberFormat = "@"

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "513030"
$q3.Cells.Item(2, 3).Value = "华安国际龙头（DAX）ETF（QDII）"
$q3.Cells.Item(2, 4).Value = "5.54"
$q3.Cells.Item(2, 5).Value = "93.57"
$q3.Cells.Item(2, 6).Value = "6.03"
$q3.Cells.Item(2, 7).Value = "0.3341"
$q3.Cells.Item(2, 8).Value = 4

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: push existing rows down one slot and insert the
#    new 2022-Q3 summary as the new row 2.
# ---------------------------------------------------------------------
$total.Range("A5:D5").Copy($total.Range("A6:D6"))
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.33

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------
# 3) Restore the originally-selected tab (the last sheet, "2020-Q4") since
#    copying/renaming sheets above shifted the active-sheet focus.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
